$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.342.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.889.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.07"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4837"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.55%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2897"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06599"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.884.53"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.90"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.156"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "87.59"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6605"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.297.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.44"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007772"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9993"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.429"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.32%  "
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.138.29"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "193.65"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.180"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.358"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.32"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.17"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -4.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.932"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.455"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.303"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09143"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.039"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05081"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -11.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.144"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7281"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01787"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.649"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9216"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.069"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.878"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.77%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.81"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4312"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.470"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.42%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.588"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +10.20%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1329"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.76"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -10.58%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.976"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05765"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -4.42%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.94"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -6.09%  "
